# Weekly update: insert a new price-record row right before the current
# row 318 ("Hortaliza, Terminal Hortofrutícola Agro Chillán - Cebolla"),
# pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 318 - shifts rows 318:360 down to 319:361
# and grows the used range to A1:R361.
$ws.Rows.Item(318).Insert()

# Populate the newly inserted row 318 with the new weekly record.
$ws.Range("A318").Value = 7
$ws.Range("B318").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C318").Value = "Ñuble"
$ws.Range("D318").Value = 44491
$ws.Range("E318").Value = 16
$ws.Range("F318").Value = 100112004
$ws.Range("G318").Value = "Cebolla"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "1a nueva(o)"
$ws.Range("J318").Value = 24000
$ws.Range("K318").Value = 1200
$ws.Range("L318").Value = 1400
$ws.Range("M318").Value = 1300
$ws.Range("N318").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O318").Value = "Región de O'Higgins"
$ws.Range("P318").Value = 130
$ws.Range("Q318").Value = 10
$ws.Range("R318").Value = "Hortaliza"
